# "Update all diagrams for Dev Guide"
#
# Content changes applied:
#   1. Rename the "AddressBook" event/handler names used throughout the
#      sequence diagram on slide 1 to "Ssenisub":
#        post(AddressBookChangedEvent)              -> post(SsenisubChangedEvent)
#        handleAddresssBookChangedEvent()            -> handleSsenisubChangedEvent()
#   2. Refresh the auto date placeholder on the slide master and every
#      slide layout from 10/16/2016 to 11/11/18.

$p = $ppt.ActivePresentation

function Replace-InShapeText {
    param(
        $shape,
        [string]$oldSub,
        [string]$newSub
    )
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    if ($full -eq $null) { return }
    $searchFrom = 0
    while ($true) {
        $idx = $full.IndexOf($oldSub, $searchFrom)
        if ($idx -lt 0) { break }
        $sub = $tr.Characters($idx + 1, $oldSub.Length)
        $sub.Text = $newSub
        $full = $tr.Text
        $searchFrom = $idx + $newSub.Length
    }
}

function Update-DateShapes {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "10/16/2016") {
                $shp.TextFrame.TextRange.Text = "11/11/18"
            }
        }
    }
}

# --- 1. Rename AddressBook* -> Ssenisub* on slide 1 -------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    Replace-InShapeText $shp "AddressBookChangedEvent" "SsenisubChangedEvent"
    Replace-InShapeText $shp "handleAddresssBookChangedEvent" "handleSsenisubChangedEvent"
}

# --- 2. Refresh the date placeholder everywhere it appears ------------
Update-DateShapes $p.SlideMaster.Shapes

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    Update-DateShapes $layout.Shapes
}
